$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's header/data rows carry one column per SAP field, e.g.
#   ... | DA: U_Regimen | DB: U_Proveedor | DC: U_Fecha_Pago(hidden) |
# "changed fecha pago": a new (visible) "U_Fecha_Pago" column is introduced
# right before the existing "U_Proveedor" column, pushing "U_Proveedor" and
# the original (hidden) "U_Fecha_Pago" column one slot to the right.
#
# Insert a blank column at DB (column 106); this shifts:
#   U_Proveedor (was DB/106, width ~13.58, visible) -> DC/107
#   U_Fecha_Pago (was DC/107, width ~12.43, hidden) -> DD/108
$ws.Columns.Item(106).Insert()

# Populate the newly inserted column (now DB/106) with the "Fecha_Pago"
# field for both the header row and the sample-data row, and give it the
# same (narrow) width as the original hidden Fecha_Pago column - this
# occurrence stays visible (not hidden).
$ws.Cells.Item(1, 106).Value = "U_Fecha_Pago"
$ws.Cells.Item(2, 106).Value = "U_Fecha_Pago"
$ws.Columns.Item(106).ColumnWidth = 11.62
